$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# New monthly rows appended after the existing data (row 361 -> 362..367)
$newRows = @(
    @(45382, 33),
    @(45412, 26),
    @(45443, 28),
    @(45473, 26),
    @(45504, 24),
    @(45535, 25)
)

$startRow = 362
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $dateSerial = $newRows[$i][0]
    $value = $newRows[$i][1]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = "m/d/yy"

    $valueCell = $ws.Cells.Item($r, 2)
    $valueCell.Value = $value
}

# Column A was resized (best-fit) to accommodate the date values
$ws.Columns.Item(1).AutoFit() | Out-Null

# Scroll the view down and leave the last new cell selected, matching
# where the author ended up after pasting the new rows in
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("B$lastRow").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 347 } catch { }
